$wb = $excel.ActiveWorkbook

# Metadata sheet: fix "Name" value and "Date" value
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B4").Value = "CompetencemetierVs"
$wsMeta.Range("B8").Value = "2025-10-29T11:46:56+00:00"

# Include sheets: swap the System URI values between the two "include" sheets
$wsInc0 = $wb.Worksheets.Item("Include #0")
$wsInc1 = $wb.Worksheets.Item("Include #1")

$wsInc0.Range("B4").Value = "https://smt.esante.gouv.fr/fhir/CodeSystem/tre-r394-competence-metier"
$wsInc1.Range("B4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R01-EnsembleSavoirFaire-CISIS/FHIR/TRE-R01-EnsembleSavoirFaire-CISIS"
